$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.596.34"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "'2.100.40"
$ws.Range("E3").Value = "  +11.05%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'248.92"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'45.36"
$ws.Range("E8").Value = "  +4.47%  "
$ws.Range("D9").Value = "'61.41"
$ws.Range("E9").Value = "  +8.21%  "
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").Value = "'0.0730"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("D12").Value = "'0.0991"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "'14.62"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "'2.404.24"
$ws.Range("E14").Value = "  +10.99%  "
$ws.Range("D15").Value = "'0.841"
$ws.Range("E15").Value = "  +6.19%  "
$ws.Range("D16").Value = "'2.098.85"
$ws.Range("E16").Value = "  +10.55%  "
$ws.Range("D17").Value = "'5.02"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "'36.655.75"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "'72.54"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'240.70"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "'12.85"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("D23").Value = "'5.02"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -9.72%  "
$ws.Range("D26").Value = "'169.99"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").Value = "'20.55"
$ws.Range("E27").Value = "  +11.49%  "
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").Value = "'2.00"
$ws.Range("E29").Value = "  -8.01%  "
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("D31").Value = "'21.99"
$ws.Range("E31").Value = "  +56.51%  "
$ws.Range("D32").Value = "'4.43"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").Value = "'0.0593"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("E34").Value = "  +15.12%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.33"
$ws.Range("E36").Value = "  +19.65%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.908"
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.05"
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("D40").Value = "'1.36"
$ws.Range("E40").Value = "  -8.85%  "
$ws.Range("D41").Value = "'1.19"
$ws.Range("E41").Value = "  +9.12%  "
$ws.Range("D42").Value = "'99.04"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("E44").Value = "  +16.31%  "
$ws.Range("D45").Value = "'16.15"
$ws.Range("E45").Value = "  -5.15%  "
$ws.Range("D46").Value = "'1.359.84"
$ws.Range("E46").Value = "  +3.28%  "
$ws.Range("D47").Value = "'0.0833"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").Value = "'2.294.04"
$ws.Range("E48").Value = "  +10.97%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("E50").Value = "  -3.73%  "
$ws.Range("D51").Value = "'3.91"
$ws.Range("E51").Value = "  +16.79%  "
